$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -7.903
$ws.Range("D4").Value = -7.763
$ws.Range("C7").Value = -13.566
$ws.Range("B8").Value = 7.086999999999999
$ws.Range("B10").Value = 5.840000000000001
$ws.Range("E10").Value = 16.694
$ws.Range("D11").Value = -7.409000000000001
$ws.Range("B12").Value = 5.691000000000001
$ws.Range("E12").Value = 17.221
$ws.Range("E13").Value = 16.43
$ws.Range("C14").Value = -13.457
$ws.Range("D14").Value = -8.202
$ws.Range("E14").Value = 16.558
$ws.Range("C15").Value = -13.745
$ws.Range("B18").Value = 5.609
$ws.Range("C18").Value = -13.483
$ws.Range("D18").Value = -8.797000000000001
$ws.Range("D19").Value = -8.255000000000001
$ws.Range("C20").Value = -12.684
$ws.Range("D21").Value = -8.220000000000001
$ws.Range("B25").Value = 6.825
$ws.Range("D27").Value = -8.073
$ws.Range("C29").Value = -11.893
$ws.Range("E29").Value = 16.808
$ws.Range("C30").Value = -12.338
$ws.Range("C31").Value = -12.198
$ws.Range("D31").Value = -7.875
$ws.Range("E32").Value = 16.448
$ws.Range("C35").Value = -12.667
$ws.Range("E35").Value = 16.494
$ws.Range("B37").Value = 8.416
$ws.Range("D38").Value = -7.892
$ws.Range("C40").Value = -12.782
$ws.Range("D42").Value = -8.300000000000001
$ws.Range("E43").Value = 16.937
$ws.Range("C44").Value = -12.395
$ws.Range("D44").Value = -7.672999999999999
$ws.Range("D47").Value = -7.683
$ws.Range("E48").Value = 17.17599999999999
$ws.Range("E49").Value = 16.349
$ws.Range("C50").Value = -12.904
$ws.Range("E50").Value = 16.435
$ws.Range("E51").Value = 16.786
$ws.Range("C54").Value = -12.448
$ws.Range("B55").Value = 5.845000000000001
$ws.Range("D56").Value = -7.994
$ws.Range("E56").Value = 16.205
$ws.Range("D58").Value = -8.134
$ws.Range("E61").Value = 16.408
$ws.Range("D65").Value = -7.831
$ws.Range("B68").Value = 5.456000000000001
$ws.Range("C68").Value = -11.307
$ws.Range("E69").Value = 17.32
$ws.Range("E71").Value = 17.208
$ws.Range("D73").Value = -8.318000000000001
$ws.Range("C76").Value = -13.46
$ws.Range("B77").Value = 5.144
$ws.Range("B78").Value = 7.723999999999999
$ws.Range("B79").Value = 5.545999999999999
$ws.Range("E79").Value = 17.523
$ws.Range("B80").Value = 8.352
$ws.Range("B81").Value = 5.875
$ws.Range("E81").Value = 16.392
$ws.Range("B82").Value = 6.218000000000001
$ws.Range("B84").Value = 6.343999999999999
$ws.Range("C87").Value = -12.772
$ws.Range("C88").Value = -12.696
$ws.Range("D90").Value = -7.555999999999999
$ws.Range("C92").Value = -11.792
$ws.Range("D92").Value = -7.472
$ws.Range("E92").Value = 16.818
$ws.Range("D94").Value = -6.953999999999999
$ws.Range("D95").Value = -7.806
$ws.Range("C96").Value = -12.705
$ws.Range("C98").Value = -13.649
$ws.Range("B101").Value = 8.975
$ws.Range("C101").Value = -13.022
$ws.Range("D101").Value = -7.784000000000001
$ws.Range("B102").Value = 7.468000000000001
$ws.Range("C102").Value = -12.906
